$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.911561666666666
$ws.Range("H2").Value = 5.734684999999999
$ws.Range("I2").Value = 0.1720155802183755
$ws.Range("J2").Value = 0.1720155802183755
$ws.Range("M2").Value = 3.483060666666667
$ws.Range("N2").Value = 10.449182
$ws.Range("O2").Value = 0.2527672867110271
$ws.Range("P2").Value = 0.2527672867110271
$ws.Range("Q2").Value = 6.658085253074444
$ws.Range("R2").Value = 59.92276727766999
$ws.Range("S2").Value = 0.0434799114838218
$ws.Range("T2").Value = 0.04347991148382179
$ws.Range("G3").Value = 1.911561666666666
$ws.Range("H3").Value = 5.734684999999999
$ws.Range("I3").Value = 0.1720155802183755
$ws.Range("J3").Value = 0.1720155802183755
$ws.Range("O3").Value = 0.3353267952677969
$ws.Range("P3").Value = 0.335326795267797
$ws.Range("Q3").Value = 8.83276637409833
$ws.Range("R3").Value = 79.49489736688497
$ws.Range("S3").Value = 0.0576814332507585
$ws.Range("T3").Value = 0.05768143325075849
$ws.Range("G4").Value = 1.911561666666666
$ws.Range("H4").Value = 5.734684999999999
$ws.Range("I4").Value = 0.1720155802183755
$ws.Range("J4").Value = 0.1720155802183755
$ws.Range("M4").Value = 2.773309666666667
$ws.Range("N4").Value = 8.319929
$ws.Range("O4").Value = 0.2012603358768551
$ws.Range("P4").Value = 0.2012603358768551
$ws.Range("Q4").Value = 5.30135244859611
$ws.Range("R4").Value = 47.71217203736499
$ws.Range("S4").Value = 0.03461991345080237
$ws.Range("T4").Value = 0.03461991345080236
$ws.Range("G5").Value = 1.911561666666666
$ws.Range("H5").Value = 5.734684999999999
$ws.Range("I5").Value = 0.1720155802183755
$ws.Range("J5").Value = 0.1720155802183755
$ws.Range("M5").Value = 2.902635666666666
$ws.Range("N5").Value = 8.707906999999999
$ws.Range("O5").Value = 0.2106455821443209
$ws.Range("P5").Value = 0.2106455821443209
$ws.Range("Q5").Value = 5.548567072699443
$ws.Range("R5").Value = 49.93710365429499
$ws.Range("S5").Value = 0.03623432203299284
$ws.Range("T5").Value = 0.03623432203299282
$ws.Range("I6").Value = 0.4009917520372743
$ws.Range("J6").Value = 0.4009917520372743
$ws.Range("M6").Value = 3.483060666666667
$ws.Range("N6").Value = 10.449182
$ws.Range("O6").Value = 0.2527672867110271
$ws.Range("P6").Value = 0.2527672867110271
$ws.Range("Q6").Value = 15.52090378938044
$ws.Range("R6").Value = 139.688134104424
$ws.Range("S6").Value = 0.1013575971559628
$ws.Range("T6").Value = 0.1013575971559628
$ws.Range("I7").Value = 0.4009917520372743
$ws.Range("J7").Value = 0.4009917520372743
$ws.Range("O7").Value = 0.3353267952677969
$ws.Range("P7").Value = 0.335326795267797
$ws.Range("S7").Value = 0.1344632791394783
$ws.Range("T7").Value = 0.1344632791394783
$ws.Range("I8").Value = 0.4009917520372743
$ws.Range("J8").Value = 0.4009917520372743
$ws.Range("M8").Value = 2.773309666666667
$ws.Range("N8").Value = 8.319929
$ws.Range("O8").Value = 0.2012603358768551
$ws.Range("P8").Value = 0.2012603358768551
$ws.Range("Q8").Value = 12.35817478760311
$ws.Range("R8").Value = 111.223573088428
$ws.Range("S8").Value = 0.08070373469887043
$ws.Range("T8").Value = 0.08070373469887043
$ws.Range("I9").Value = 0.4009917520372743
$ws.Range("J9").Value = 0.4009917520372743
$ws.Range("M9").Value = 2.902635666666666
$ws.Range("N9").Value = 8.707906999999999
$ws.Range("O9").Value = 0.2106455821443209
$ws.Range("P9").Value = 0.2106455821443209
$ws.Range("Q9").Value = 12.93446575568044
$ws.Range("R9").Value = 116.410191801124
$ws.Range("S9").Value = 0.08446714104296282
$ws.Range("T9").Value = 0.08446714104296281
$ws.Range("G10").Value = 4.603447666666667
$ws.Range("H10").Value = 13.810343
$ws.Range("I10").Value = 0.4142501574471451
$ws.Range("J10").Value = 0.4142501574471449
$ws.Range("M10").Value = 3.483060666666667
$ws.Range("N10").Value = 10.449182
$ws.Range("O10").Value = 0.2527672867110271
$ws.Range("P10").Value = 0.2527672867110271
$ws.Range("Q10").Value = 16.03408749882511
$ws.Range("R10").Value = 144.306787489426
$ws.Range("S10").Value = 0.1047088883175306
$ws.Range("T10").Value = 0.1047088883175306
$ws.Range("G11").Value = 4.603447666666667
$ws.Range("H11").Value = 13.810343
$ws.Range("I11").Value = 0.4142501574471451
$ws.Range("J11").Value = 0.4142501574471449
$ws.Range("O11").Value = 0.3353267952677969
$ws.Range("P11").Value = 0.335326795267797
$ws.Range("Q11").Value = 21.27118285750033
$ws.Range("R11").Value = 191.440645717503
$ws.Range("S11").Value = 0.1389091777359314
$ws.Range("T11").Value = 0.1389091777359314
$ws.Range("G12").Value = 4.603447666666667
$ws.Range("H12").Value = 13.810343
$ws.Range("I12").Value = 0.4142501574471451
$ws.Range("J12").Value = 0.4142501574471449
$ws.Range("M12").Value = 2.773309666666667
$ws.Range("N12").Value = 8.319929
$ws.Range("O12").Value = 0.2012603358768551
$ws.Range("P12").Value = 0.2012603358768551
$ws.Range("Q12").Value = 12.76678591396078
$ws.Range("R12").Value = 114.901073225647
$ws.Range("S12").Value = 0.08337212582485252
$ws.Range("T12").Value = 0.08337212582485251
$ws.Range("G13").Value = 4.603447666666667
$ws.Range("H13").Value = 13.810343
$ws.Range("I13").Value = 0.4142501574471451
$ws.Range("J13").Value = 0.4142501574471449
$ws.Range("M13").Value = 2.902635666666666
$ws.Range("N13").Value = 8.707906999999999
$ws.Range("O13").Value = 0.2106455821443209
$ws.Range("P13").Value = 0.2106455821443209
$ws.Range("Q13").Value = 13.36213138690011
$ws.Range("R13").Value = 120.259182482101
$ws.Range("S13").Value = 0.08725996556883045
$ws.Range("T13").Value = 0.08725996556883044
$ws.Range("G14").Value = 0.141604
$ws.Range("H14").Value = 0.424812
$ws.Range("I14").Value = 0.01274251029720526
$ws.Range("J14").Value = 0.01274251029720526
$ws.Range("M14").Value = 3.483060666666667
$ws.Range("N14").Value = 10.449182
$ws.Range("O14").Value = 0.2527672867110271
$ws.Range("P14").Value = 0.2527672867110271
$ws.Range("Q14").Value = 0.4932153226426667
$ws.Range("R14").Value = 4.438937903784001
$ws.Range("S14").Value = 0.003220889753711896
$ws.Range("T14").Value = 0.003220889753711897
$ws.Range("G15").Value = 0.141604
$ws.Range("H15").Value = 0.424812
$ws.Range("I15").Value = 0.01274251029720526
$ws.Range("J15").Value = 0.01274251029720526
$ws.Range("O15").Value = 0.3353267952677969
$ws.Range("P15").Value = 0.335326795267797
$ws.Range("Q15").Value = 0.654310594028
$ws.Range("R15").Value = 5.888795346252
$ws.Range("S15").Value = 0.004272905141628742
$ws.Range("T15").Value = 0.004272905141628742
$ws.Range("G16").Value = 0.141604
$ws.Range("H16").Value = 0.424812
$ws.Range("I16").Value = 0.01274251029720526
$ws.Range("J16").Value = 0.01274251029720526
$ws.Range("M16").Value = 2.773309666666667
$ws.Range("N16").Value = 8.319929
$ws.Range("O16").Value = 0.2012603358768551
$ws.Range("P16").Value = 0.2012603358768551
$ws.Range("Q16").Value = 0.3927117420386667
$ws.Range("R16").Value = 3.534405678348
$ws.Range("S16").Value = 0.002564561902329815
$ws.Range("T16").Value = 0.002564561902329815
$ws.Range("G17").Value = 0.141604
$ws.Range("H17").Value = 0.424812
$ws.Range("I17").Value = 0.01274251029720526
$ws.Range("J17").Value = 0.01274251029720526
$ws.Range("M17").Value = 2.902635666666666
$ws.Range("N17").Value = 8.707906999999999
$ws.Range("O17").Value = 0.2106455821443209
$ws.Range("P17").Value = 0.2106455821443209
$ws.Range("Q17").Value = 0.4110248209426666
$ws.Range("R17").Value = 3.699223388484
$ws.Range("S17").Value = 0.002684153499534805
$ws.Range("T17").Value = 0.002684153499534805
